$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing header "Category" to A1, matching the formatting of the
# other header cells in row 1 (bold, bordered, centered - style index 1).
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The category cells A2:A46 previously used the header style; strip that
# formatting so they match the plain (unstyled) data cells, e.g. B2.
$ws.Range("B2").Copy()
$ws.Range("A2:A46").PasteSpecial(-4122)
